$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename header columns on existing sheets
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "PO Forecast"

# Match sheet-level formatting (outline properties) used by the other sheets
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Match page margins used by the other sheets (values are in points; 72pt = 1in)
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Reuse exact header/date cell formatting from "Weekly Quantity"
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A63").PasteSpecial(-4122)

# Forecast data rows 2-63
$newSheet.Cells.Item(2,1).Value = 44990.99999999999
$newSheet.Cells.Item(2,2).Value = 41
$newSheet.Cells.Item(2,3).Value = -89.50606858047293
$newSheet.Cells.Item(2,4).Value = 179.7957006340748
$newSheet.Cells.Item(3,1).Value = 44997.99999999999
$newSheet.Cells.Item(3,2).Value = 41
$newSheet.Cells.Item(3,3).Value = -102.1342346969065
$newSheet.Cells.Item(3,4).Value = 172.7035634709087
$newSheet.Cells.Item(4,1).Value = 45011.99999999999
$newSheet.Cells.Item(4,2).Value = 43
$newSheet.Cells.Item(4,3).Value = -93.45419251997258
$newSheet.Cells.Item(4,4).Value = 172.7342859123707
$newSheet.Cells.Item(5,1).Value = 45025.99999999999
$newSheet.Cells.Item(5,2).Value = 44
$newSheet.Cells.Item(5,3).Value = -95.37181522762005
$newSheet.Cells.Item(5,4).Value = 173.4782025932979
$newSheet.Cells.Item(6,1).Value = 45039.99999999999
$newSheet.Cells.Item(6,2).Value = 46
$newSheet.Cells.Item(6,3).Value = -83.87047367189332
$newSheet.Cells.Item(6,4).Value = 183.0561959537235
$newSheet.Cells.Item(7,1).Value = 45046.99999999999
$newSheet.Cells.Item(7,2).Value = 46
$newSheet.Cells.Item(7,3).Value = -86.0565346708544
$newSheet.Cells.Item(7,4).Value = 172.0815755254318
$newSheet.Cells.Item(8,1).Value = 45053.99999999999
$newSheet.Cells.Item(8,2).Value = 47
$newSheet.Cells.Item(8,3).Value = -85.08141748213735
$newSheet.Cells.Item(8,4).Value = 186.892664825666
$newSheet.Cells.Item(9,1).Value = 45060.99999999999
$newSheet.Cells.Item(9,2).Value = 48
$newSheet.Cells.Item(9,3).Value = -83.16013135943128
$newSheet.Cells.Item(9,4).Value = 178.3754046558724
$newSheet.Cells.Item(10,1).Value = 45067.99999999999
$newSheet.Cells.Item(10,2).Value = 49
$newSheet.Cells.Item(10,3).Value = -83.85073579688374
$newSheet.Cells.Item(10,4).Value = 181.8080493927991
$newSheet.Cells.Item(11,1).Value = 45074.99999999999
$newSheet.Cells.Item(11,2).Value = 49
$newSheet.Cells.Item(11,3).Value = -77.04770440121435
$newSheet.Cells.Item(11,4).Value = 180.6361588223857
$newSheet.Cells.Item(12,1).Value = 45081.99999999999
$newSheet.Cells.Item(12,2).Value = 50
$newSheet.Cells.Item(12,3).Value = -77.31079287670971
$newSheet.Cells.Item(12,4).Value = 179.1027588761264
$newSheet.Cells.Item(13,1).Value = 45088.99999999999
$newSheet.Cells.Item(13,2).Value = 51
$newSheet.Cells.Item(13,3).Value = -86.1287455283346
$newSheet.Cells.Item(13,4).Value = 191.295203910469
$newSheet.Cells.Item(14,1).Value = 45095.99999999999
$newSheet.Cells.Item(14,2).Value = 52
$newSheet.Cells.Item(14,3).Value = -81.97446575108042
$newSheet.Cells.Item(14,4).Value = 185.9687576311161
$newSheet.Cells.Item(15,1).Value = 45102.99999999999
$newSheet.Cells.Item(15,2).Value = 52
$newSheet.Cells.Item(15,3).Value = -80.82025799688719
$newSheet.Cells.Item(15,4).Value = 183.5826674569251
$newSheet.Cells.Item(16,1).Value = 45109.99999999999
$newSheet.Cells.Item(16,2).Value = 53
$newSheet.Cells.Item(16,3).Value = -94.42882673831336
$newSheet.Cells.Item(16,4).Value = 179.4854933884369
$newSheet.Cells.Item(17,1).Value = 45116.99999999999
$newSheet.Cells.Item(17,2).Value = 54
$newSheet.Cells.Item(17,3).Value = -77.02516326182557
$newSheet.Cells.Item(17,4).Value = 187.1012605418837
$newSheet.Cells.Item(18,1).Value = 45130.99999999999
$newSheet.Cells.Item(18,2).Value = 55
$newSheet.Cells.Item(18,3).Value = -78.18690059249755
$newSheet.Cells.Item(18,4).Value = 185.8352624544149
$newSheet.Cells.Item(19,1).Value = 45137.99999999999
$newSheet.Cells.Item(19,2).Value = 56
$newSheet.Cells.Item(19,3).Value = -69.15390274680288
$newSheet.Cells.Item(19,4).Value = 192.682248830942
$newSheet.Cells.Item(20,1).Value = 45165.99999999999
$newSheet.Cells.Item(20,2).Value = 59
$newSheet.Cells.Item(20,3).Value = -77.38715940872845
$newSheet.Cells.Item(20,4).Value = 178.3257491138096
$newSheet.Cells.Item(21,1).Value = 45179.99999999999
$newSheet.Cells.Item(21,2).Value = 60
$newSheet.Cells.Item(21,3).Value = -72.86308029547567
$newSheet.Cells.Item(21,4).Value = 194.751379990795
$newSheet.Cells.Item(22,1).Value = 45186.99999999999
$newSheet.Cells.Item(22,2).Value = 61
$newSheet.Cells.Item(22,3).Value = -73.61500877339516
$newSheet.Cells.Item(22,4).Value = 191.8454596156341
$newSheet.Cells.Item(23,1).Value = 45193.99999999999
$newSheet.Cells.Item(23,2).Value = 62
$newSheet.Cells.Item(23,3).Value = -68.55901544253433
$newSheet.Cells.Item(23,4).Value = 200.4273245778207
$newSheet.Cells.Item(24,1).Value = 45214.99999999999
$newSheet.Cells.Item(24,2).Value = 64
$newSheet.Cells.Item(24,3).Value = -74.16234316939151
$newSheet.Cells.Item(24,4).Value = 194.6041687462691
$newSheet.Cells.Item(25,1).Value = 45221.99999999999
$newSheet.Cells.Item(25,2).Value = 64
$newSheet.Cells.Item(25,3).Value = -68.39610949551252
$newSheet.Cells.Item(25,4).Value = 201.8970993431581
$newSheet.Cells.Item(26,1).Value = 45235.99999999999
$newSheet.Cells.Item(26,2).Value = 66
$newSheet.Cells.Item(26,3).Value = -67.90775758380006
$newSheet.Cells.Item(26,4).Value = 201.6914081201685
$newSheet.Cells.Item(27,1).Value = 45242.99999999999
$newSheet.Cells.Item(27,2).Value = 67
$newSheet.Cells.Item(27,3).Value = -71.18075566966778
$newSheet.Cells.Item(27,4).Value = 195.0971862846638
$newSheet.Cells.Item(28,1).Value = 45249.99999999999
$newSheet.Cells.Item(28,2).Value = 67
$newSheet.Cells.Item(28,3).Value = -65.14305839833482
$newSheet.Cells.Item(28,4).Value = 195.7171030589973
$newSheet.Cells.Item(29,1).Value = 45256.99999999999
$newSheet.Cells.Item(29,2).Value = 68
$newSheet.Cells.Item(29,3).Value = -72.17878921177009
$newSheet.Cells.Item(29,4).Value = 202.387267402619
$newSheet.Cells.Item(30,1).Value = 45270.99999999999
$newSheet.Cells.Item(30,2).Value = 69
$newSheet.Cells.Item(30,3).Value = -73.34111808371118
$newSheet.Cells.Item(30,4).Value = 198.4290479865792
$newSheet.Cells.Item(31,1).Value = 45277.99999999999
$newSheet.Cells.Item(31,2).Value = 70
$newSheet.Cells.Item(31,3).Value = -56.67577420658425
$newSheet.Cells.Item(31,4).Value = 215.0835933481424
$newSheet.Cells.Item(32,1).Value = 45298.99999999999
$newSheet.Cells.Item(32,2).Value = 72
$newSheet.Cells.Item(32,3).Value = -59.58270193933155
$newSheet.Cells.Item(32,4).Value = 213.2818235495887
$newSheet.Cells.Item(33,1).Value = 45312.99999999999
$newSheet.Cells.Item(33,2).Value = 74
$newSheet.Cells.Item(33,3).Value = -54.60455913238223
$newSheet.Cells.Item(33,4).Value = 209.2890730237712
$newSheet.Cells.Item(34,1).Value = 45326.99999999999
$newSheet.Cells.Item(34,2).Value = 75
$newSheet.Cells.Item(34,3).Value = -51.06743726399991
$newSheet.Cells.Item(34,4).Value = 212.1685816866676
$newSheet.Cells.Item(35,1).Value = 45333.99999999999
$newSheet.Cells.Item(35,2).Value = 76
$newSheet.Cells.Item(35,3).Value = -60.87863763053552
$newSheet.Cells.Item(35,4).Value = 198.8274402744682
$newSheet.Cells.Item(36,1).Value = 45361.99999999999
$newSheet.Cells.Item(36,2).Value = 79
$newSheet.Cells.Item(36,3).Value = -59.61924189774241
$newSheet.Cells.Item(36,4).Value = 198.9335088402317
$newSheet.Cells.Item(37,1).Value = 45368.99999999999
$newSheet.Cells.Item(37,2).Value = 79
$newSheet.Cells.Item(37,3).Value = -56.84614039383105
$newSheet.Cells.Item(37,4).Value = 206.2487291476089
$newSheet.Cells.Item(38,1).Value = 45375.99999999999
$newSheet.Cells.Item(38,2).Value = 80
$newSheet.Cells.Item(38,3).Value = -55.13084104117181
$newSheet.Cells.Item(38,4).Value = 212.9843104681838
$newSheet.Cells.Item(39,1).Value = 45382.99999999999
$newSheet.Cells.Item(39,2).Value = 81
$newSheet.Cells.Item(39,3).Value = -50.76954443525399
$newSheet.Cells.Item(39,4).Value = 204.5016200197873
$newSheet.Cells.Item(40,1).Value = 45396.99999999999
$newSheet.Cells.Item(40,2).Value = 82
$newSheet.Cells.Item(40,3).Value = -59.58941447617529
$newSheet.Cells.Item(40,4).Value = 209.2721255467738
$newSheet.Cells.Item(41,1).Value = 45403.99999999999
$newSheet.Cells.Item(41,2).Value = 83
$newSheet.Cells.Item(41,3).Value = -50.35887676727876
$newSheet.Cells.Item(41,4).Value = 223.163167085258
$newSheet.Cells.Item(42,1).Value = 45410.99999999999
$newSheet.Cells.Item(42,2).Value = 84
$newSheet.Cells.Item(42,3).Value = -38.51324335637922
$newSheet.Cells.Item(42,4).Value = 215.3134641466439
$newSheet.Cells.Item(43,1).Value = 45417.99999999999
$newSheet.Cells.Item(43,2).Value = 84
$newSheet.Cells.Item(43,3).Value = -45.65863728436909
$newSheet.Cells.Item(43,4).Value = 217.2923677341927
$newSheet.Cells.Item(44,1).Value = 45424.99999999999
$newSheet.Cells.Item(44,2).Value = 85
$newSheet.Cells.Item(44,3).Value = -44.37258817790311
$newSheet.Cells.Item(44,4).Value = 220.8570796239159
$newSheet.Cells.Item(45,1).Value = 45459.99999999999
$newSheet.Cells.Item(45,2).Value = 89
$newSheet.Cells.Item(45,3).Value = -44.60167700595971
$newSheet.Cells.Item(45,4).Value = 211.7279555178088
$newSheet.Cells.Item(46,1).Value = 45466.99999999999
$newSheet.Cells.Item(46,2).Value = 89
$newSheet.Cells.Item(46,3).Value = -33.77454135009368
$newSheet.Cells.Item(46,4).Value = 232.8238515763462
$newSheet.Cells.Item(47,1).Value = 45480.99999999999
$newSheet.Cells.Item(47,2).Value = 91
$newSheet.Cells.Item(47,3).Value = -43.13864414665338
$newSheet.Cells.Item(47,4).Value = 219.3779010775385
$newSheet.Cells.Item(48,1).Value = 45487.99999999999
$newSheet.Cells.Item(48,2).Value = 92
$newSheet.Cells.Item(48,3).Value = -42.12111442780365
$newSheet.Cells.Item(48,4).Value = 225.6580613606374
$newSheet.Cells.Item(49,1).Value = 45494.99999999999
$newSheet.Cells.Item(49,2).Value = 92
$newSheet.Cells.Item(49,3).Value = -37.68556528628482
$newSheet.Cells.Item(49,4).Value = 222.9284515360052
$newSheet.Cells.Item(50,1).Value = 45501.99999999999
$newSheet.Cells.Item(50,2).Value = 93
$newSheet.Cells.Item(50,3).Value = -39.71769902240941
$newSheet.Cells.Item(50,4).Value = 224.4124304749129
$newSheet.Cells.Item(51,1).Value = 45508.99999999999
$newSheet.Cells.Item(51,2).Value = 94
$newSheet.Cells.Item(51,3).Value = -34.16100298626581
$newSheet.Cells.Item(51,4).Value = 231.9076509029759
$newSheet.Cells.Item(52,1).Value = 45515.99999999999
$newSheet.Cells.Item(52,2).Value = 94
$newSheet.Cells.Item(52,3).Value = -33.9444576702701
$newSheet.Cells.Item(52,4).Value = 229.2046116642845
$newSheet.Cells.Item(53,1).Value = 45522.99999999999
$newSheet.Cells.Item(53,2).Value = 95
$newSheet.Cells.Item(53,3).Value = -37.88884242486667
$newSheet.Cells.Item(53,4).Value = 239.3693566359052
$newSheet.Cells.Item(54,1).Value = 45529.99999999999
$newSheet.Cells.Item(54,2).Value = 96
$newSheet.Cells.Item(54,3).Value = -34.65896578218703
$newSheet.Cells.Item(54,4).Value = 233.6175535044598
$newSheet.Cells.Item(55,1).Value = 45536.99999999999
$newSheet.Cells.Item(55,2).Value = 97
$newSheet.Cells.Item(55,3).Value = -37.89812821312167
$newSheet.Cells.Item(55,4).Value = 228.1334482286829
$newSheet.Cells.Item(56,1).Value = 45543.99999999999
$newSheet.Cells.Item(56,2).Value = 97
$newSheet.Cells.Item(56,3).Value = -36.47231347910745
$newSheet.Cells.Item(56,4).Value = 231.2108266560249
$newSheet.Cells.Item(57,1).Value = 45550.99999999999
$newSheet.Cells.Item(57,2).Value = 98
$newSheet.Cells.Item(57,3).Value = -29.5925225822459
$newSheet.Cells.Item(57,4).Value = 231.3198806285499
$newSheet.Cells.Item(58,1).Value = 45557.99999999999
$newSheet.Cells.Item(58,2).Value = 99
$newSheet.Cells.Item(58,3).Value = -41.37707602714326
$newSheet.Cells.Item(58,4).Value = 226.6649817629838
$newSheet.Cells.Item(59,1).Value = 45564.99999999999
$newSheet.Cells.Item(59,2).Value = 99
$newSheet.Cells.Item(59,3).Value = -23.13701889821078
$newSheet.Cells.Item(59,4).Value = 233.8623324681299
$newSheet.Cells.Item(60,1).Value = 45571.99999999999
$newSheet.Cells.Item(60,2).Value = 100
$newSheet.Cells.Item(60,3).Value = -24.35929472797083
$newSheet.Cells.Item(60,4).Value = 224.9723641279391
$newSheet.Cells.Item(61,1).Value = 45578.99999999999
$newSheet.Cells.Item(61,2).Value = 101
$newSheet.Cells.Item(61,3).Value = -29.44259352608129
$newSheet.Cells.Item(61,4).Value = 231.3619123069622
$newSheet.Cells.Item(62,1).Value = 45585.99999999999
$newSheet.Cells.Item(62,2).Value = 102
$newSheet.Cells.Item(62,3).Value = -27.32735210404779
$newSheet.Cells.Item(62,4).Value = 236.9718910869815
$newSheet.Cells.Item(63,1).Value = 45592.99999999999
$newSheet.Cells.Item(63,2).Value = 102
$newSheet.Cells.Item(63,3).Value = -30.92341360175019
$newSheet.Cells.Item(63,4).Value = 222.7479214096296

# Restore original active sheet/tab selection
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
